$d = $word.ActiveDocument
$cr = [char]13

# The document currently ends with an empty "ListBullet" paragraph (right
# before the final section break) that follows the bold "NGINX" heading.
# Add the NGINX script content as three new paragraphs after it: a
# descriptive line, a blank spacer line, and a second descriptive line -
# each inheriting the same ListBullet paragraph formatting.

$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$r = $lastPara.Range
$r.Collapse(0)
$r.Text = $cr + "NGINX is a tool that is used for; web serving, reverse proxying, caching and load balancing."

$n = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($n)
$r2 = $p.Range
$r2.Collapse(0)
$r2.Text = $cr

$n = $d.Paragraphs.Count
$p2 = $d.Paragraphs.Item($n)
$r3 = $p2.Range
$r3.Collapse(0)
$r3.Text = $cr + "We have used it as a load balancer in order to access our different docker containers in our swarm. It automatically sends the user to an available container."

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
